$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "E3"   = 16.392
    "C7"   = -12.675
    "A8"   = -22.188
    "A10"  = -21.736
    "A12"  = -21.683
    "C15"  = -13.636
    "A18"  = -21.572
    "C18"  = -10.643
    "E18"  = 17.793
    "E19"  = 16.538
    "C20"  = -12.183
    "E27"  = 16.501
    "C29"  = -11.987
    "C30"  = -13.347
    "C31"  = -13.519
    "E31"  = 16.238
    "A37"  = -20.029
    "E38"  = 16.7
    "C40"  = -12.782
    "E42"  = 16.576
    "E44"  = 16.682
    "E47"  = 16.32
    "C50"  = -13.371
    "A55"  = -21.868
    "E58"  = 16.602
    "E65"  = 17.301
    "A68"  = -21.736
    "C68"  = -11.001
    "E73"  = 16.546
    "C76"  = -12.72
    "A77"  = -20.843
    "A78"  = -20.134
    "A81"  = -21.869
    "A82"  = -22.152
    "C87"  = -13.199
    "C88"  = -13.091
    "E90"  = 16.437
    "E94"  = 17.828
    "E95"  = 17.399
    "C96"  = -12.81
    "C98"  = -13.201
    "C101" = -12.747
    "E101" = 16.701
    "C102" = -13.091
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
